$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), matching the style of existing headers (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-25: column I = 1 (constant), column J = same as column H
for ($r = 2; $r -le 25; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
